# Apply the diff to the document:
#  1. Extend the first INSERT ... VALUES list with four more ('Kreuz ...') tuples.
#  2. Add a "straße" character varying(4) column to the "routen" table definition.
#  3. Lower-case the second "Insert into" -> "insert into".
#  4. Add the straße column + new rows to the second INSERT ... VALUES list.
#  5. Move the hidden _GoBack bookmark from the end of the first paragraph to the
#     last (second empty) paragraph.

$d = $word.ActiveDocument

# --- 1: extend the city-list statement with the four new "Kreuz" tuples ---
$old1 = "('Hildesheim'),('Cottbus');"
$new1 = "('Hildesheim'),('Cottbus'), ('Kreuz A13/A15'), ('Kreuz A19/A24'),('Kreuz A10/A24'), ('Kreuz A2/A10');"
$r1 = $d.Content
$r1.Find.Execute($old1) | Out-Null
$r1.Text = $new1

# --- 2: add the straße column to the CREATE TABLE routen(...) statement ---
$old2 = "strecke integer NOT NULL,PRIMARY KEY (start, "
$new2 = 'strecke integer NOT NULL,"straße" character varying(4) NOT NULL,PRIMARY KEY (start, '
$r2 = $d.Content
$r2.Find.Execute($old2) | Out-Null
$r2.Text = $new2

# --- 3: lower-case the second "Insert into" ---
$old3 = "OIDS = FALSE);Insert into "
$new3 = "OIDS = FALSE);insert into "
$r3 = $d.Content
$r3.Find.Execute($old3) | Out-Null
$r3.Text = $new3

# --- 4: add straße column + replace the values list on the routen insert ---
$old4 = " (start,ziel,strecke) values (6,61,83),(61,23,73),(23,3,50),(35,2,51);"
$new4 = " (start,ziel,strecke,straße) values (1,44,21,'A115'),(1,81,63,'A13'),(80,81,31,'A15'),(1,83,20,'A10'),(82,83,64,'A24'),(39,82,114,'A19'),(2,83,237,'A24'),(44,84,20,'A10'),(32,84,84,'A2'),(10,44,127,'A9');"
$r4 = $d.Content
$r4.Find.Execute($old4) | Out-Null
$r4.Text = $new4

# --- 5: relocate the hidden _GoBack bookmark to the last paragraph ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$lastPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
